$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.872.63"
$ws.Range("E2").Value = "'  +2.85%  "
$ws.Range("D3").Value = "'2.665.20"
$ws.Range("E3").Value = "'  +2.85%  "
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'600.99"
$ws.Range("E5").Value = "'  +2.20%  "
$ws.Range("D6").Value = "'155.80"
$ws.Range("E6").Value = "'  +4.63%  "
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("D8").Value = "'0.547"
$ws.Range("E8").Value = "'  +0.74%  "
$ws.Range("D9").Value = "'2.665.03"
$ws.Range("E9").Value = "'  +2.92%  "
$ws.Range("E10").Value = "'  +13.76%  "
$ws.Range("D12").Value = "'5.24"
$ws.Range("E12").Value = "'  +1.87%  "
$ws.Range("D13").Value = "'0.350"
$ws.Range("E13").Value = "'  +2.27%  "
$ws.Range("D14").Value = "'28.07"
$ws.Range("E14").Value = "'  +3.82%  "
$ws.Range("E15").Value = "'  +6.32%  "
$ws.Range("D16").Value = "'3.148.72"
$ws.Range("E16").Value = "'  +2.82%  "
$ws.Range("D17").Value = "'68.740.42"
$ws.Range("E17").Value = "'  +2.61%  "
$ws.Range("D18").Value = "'2.672.39"
$ws.Range("E18").Value = "'  +3.14%  "
$ws.Range("D19").Value = "'11.47"
$ws.Range("E19").Value = "'  +4.75%  "
$ws.Range("D20").Value = "'367.23"
$ws.Range("E20").Value = "'  +1.64%  "
$ws.Range("E21").Value = "'  +2.35%  "
$ws.Range("E22").Value = "'  +0.12%  "
$ws.Range("D23").Value = "'4.89"
$ws.Range("E23").Value = "'  +1.67%  "
$ws.Range("E24").Value = "'  +5.47%  "
$ws.Range("E25").Value = "'  +1.05%  "
$ws.Range("E26").Value = "'  +0.04%  "
$ws.Range("D27").Value = "'10.06"
$ws.Range("E27").Value = "'  +1.45%  "
$ws.Range("D28").Value = "'0.0000106"
$ws.Range("E28").Value = "'  +9.04%  "
$ws.Range("D29").Value = "'2.803.36"
$ws.Range("E29").Value = "'  +3.38%  "
$ws.Range("E30").Value = "'  -0.10%  "
$ws.Range("D31").Value = "'578.68"
$ws.Range("E31").Value = "'  +0.87%  "
$ws.Range("E32").Value = "'  +4.86%  "
$ws.Range("D33").Value = "'8.00"
$ws.Range("E33").Value = "'  +5.81%  "
$ws.Range("D34").Value = "'1.86"
$ws.Range("E34").Value = "'  +3.85%  "
$ws.Range("E35").Value = "'  +5.76%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "'  +0.05%  "
$ws.Range("D37").Value = "'1.54"
$ws.Range("E37").Value = "'  +4.27%  "
$ws.Range("D38").Value = "'159.45"
$ws.Range("E38").Value = "'  +2.03%  "
$ws.Range("E39").Value = "'  +5.36%  "
$ws.Range("D40").Value = "'19.32"
$ws.Range("E40").Value = "'  +2.55%  "
$ws.Range("D41").Value = "'5.43"
$ws.Range("E41").Value = "'  +5.41%  "
$ws.Range("D42").Value = "'0.369"
$ws.Range("E42").Value = "'  +1.34%  "
$ws.Range("D43").Value = "'2.67"
$ws.Range("E43").Value = "'  +7.76%  "
$ws.Range("D44").Value = "'17.75"
$ws.Range("E44").Value = "'  +5.92%  "
$ws.Range("E45").Value = "'  +14.41%  "
$ws.Range("D46").Value = "'40.66"
$ws.Range("E46").Value = "'  -0.26%  "
$ws.Range("E47").Value = "'  +0.12%  "
$ws.Range("D48").Value = "'156.71"
$ws.Range("E48").Value = "'  +3.43%  "
$ws.Range("D49").Value = "'3.75"
$ws.Range("E49").Value = "'  +1.31%  "
$ws.Range("E50").Value = "'  +3.21%  "
$ws.Range("D51").Value = "'22.08"
$ws.Range("E51").Value = "'  +4.29%  "
